# Auto-generated edit script applying numeric updates to Sagittarius_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 316.27777
$ws.Range("I39").Value = 339.2857
$ws.Range("J39").Value = 235.75
$ws.Range("K39").Value = 1017.8571
$ws.Range("L39").Value = 707.25
$ws.Range("M39").Value = -721.8571000000001
$ws.Range("N39").Value = -1299.25
$ws.Range("H51").Value = 9999
$ws.Range("J51").Value = 9999
$ws.Range("L51").Value = 9999
$ws.Range("N51").Value = -10967
$ws.Range("H127").Value = 2580.6667
$ws.Range("I127").Value = 2658.077
$ws.Range("J127").Value = 2379.4
$ws.Range("K127").Value = 7974.231000000001
$ws.Range("L127").Value = 7138.200000000001
$ws.Range("M127").Value = -3014.231000000001
$ws.Range("N127").Value = -17058.2
$ws.Range("H135").Value = 1221.9445
$ws.Range("I135").Value = 1291.0714
$ws.Range("J135").Value = 980
$ws.Range("K135").Value = 11619.6426
$ws.Range("L135").Value = 8820
$ws.Range("M135").Value = -9084.642600000001
$ws.Range("N135").Value = -13890

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1801.2
$ws.Range("I2").Value = 652.5714
$ws.Range("K2").Value = 652.5714
$ws.Range("M2").Value = -539.5714
$ws.Range("H32").Value = 4056109.5
$ws.Range("I32").Value = 4121491
$ws.Range("K32").Value = 4121491
$ws.Range("M32").Value = -4121204
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H61").Value = 10850
$ws.Range("I61").Value = 3500
$ws.Range("K61").Value = 3500
$ws.Range("M61").Value = -3288
$ws.Range("H74").Value = 752
$ws.Range("I74").Value = 683.4
$ws.Range("J74").Value = 866.3333
$ws.Range("K74").Value = 683.4
$ws.Range("L74").Value = 866.3333
$ws.Range("M74").Value = 190.6
$ws.Range("N74").Value = -2614.3333
$ws.Range("H77").Value = 752
$ws.Range("I77").Value = 683.4
$ws.Range("J77").Value = 866.3333
$ws.Range("K77").Value = 3417
$ws.Range("L77").Value = 4331.6665
$ws.Range("M77").Value = 951
$ws.Range("N77").Value = -13067.6665
$ws.Range("H116").Value = 1801.2
$ws.Range("I116").Value = 652.5714
$ws.Range("K116").Value = 652.5714
$ws.Range("M116").Value = 1641.4286
$ws.Range("H122").Value = 25799.77
$ws.Range("I122").Value = 29881.545
$ws.Range("J122").Value = 3350
$ws.Range("K122").Value = 89644.635
$ws.Range("L122").Value = 10050
$ws.Range("M122").Value = -87194.635
$ws.Range("N122").Value = -14950
$ws.Range("H132").Value = 1671.909
$ws.Range("I132").Value = 1599.5555
$ws.Range("K132").Value = 4798.666499999999
$ws.Range("M132").Value = -2268.666499999999
$ws.Range("H136").Value = 10850
$ws.Range("I136").Value = 3500
$ws.Range("K136").Value = 10500
$ws.Range("M136").Value = -7950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1801.2
$ws.Range("I3").Value = 652.5714
$ws.Range("K3").Value = 652.5714
$ws.Range("M3").Value = -538.5714
$ws.Range("H20").Value = 1884.9166
$ws.Range("I20").Value = 2016.5
$ws.Range("J20").Value = 1621.75
$ws.Range("K20").Value = 2016.5
$ws.Range("L20").Value = 1621.75
$ws.Range("M20").Value = -1769.5
$ws.Range("N20").Value = -2115.75
$ws.Range("H107").Value = 3293.9333
$ws.Range("I107").Value = 3058.0908
$ws.Range("K107").Value = 3058.0908
$ws.Range("M107").Value = -1138.0908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 678.4286
$ws.Range("I122").Value = 573.75
$ws.Range("J122").Value = 818
$ws.Range("K122").Value = 1721.25
$ws.Range("L122").Value = 2454
$ws.Range("M122").Value = 728.75
$ws.Range("N122").Value = -7354
$ws.Range("H132").Value = 3705.625
$ws.Range("I132").Value = 3789.6667
$ws.Range("J132").Value = 2445
$ws.Range("K132").Value = 11369.0001
$ws.Range("L132").Value = 7335
$ws.Range("M132").Value = -8839.000100000001
$ws.Range("N132").Value = -12395
$ws.Range("H134").Value = 2543.6667
$ws.Range("I134").Value = 2543.6667
$ws.Range("K134").Value = 7631.000100000001
$ws.Range("M134").Value = -5096.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 274.14285
$ws.Range("I6").Value = 274.14285
$ws.Range("K6").Value = 822.4285500000001
$ws.Range("M6").Value = -709.4285500000001
$ws.Range("H40").Value = 31.6
$ws.Range("I40").Value = 46.333332
$ws.Range("J40").Value = 9.5
$ws.Range("K40").Value = 185.333328
$ws.Range("L40").Value = 38
$ws.Range("M40").Value = -116.333328
$ws.Range("N40").Value = -176
$ws.Range("H69").Value = 1832
$ws.Range("I69").Value = 650
$ws.Range("K69").Value = 1950
$ws.Range("M69").Value = -1139
$ws.Range("H72").Value = 1832
$ws.Range("I72").Value = 650
$ws.Range("K72").Value = 5850
$ws.Range("M72").Value = -1794
$ws.Range("H81").Value = 2360
$ws.Range("H84").Value = 2360
$ws.Range("H114").Value = 2497.5
$ws.Range("I114").Value = 2514
$ws.Range("J114").Value = 2481
$ws.Range("K114").Value = 7542
$ws.Range("L114").Value = 7443
$ws.Range("M114").Value = -4288
$ws.Range("N114").Value = -13951

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1049.5
$ws.Range("I97").Value = 1099.3334
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 1099.3334
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -603.3334
$ws.Range("N97").Value = -1892

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1119.0769
$ws.Range("I22").Value = 690.2
$ws.Range("K22").Value = 690.2
$ws.Range("M22").Value = -395.2
$ws.Range("H27").Value = 1119.0769
$ws.Range("I27").Value = 690.2
$ws.Range("K27").Value = 690.2
$ws.Range("M27").Value = -583.2
$ws.Range("H40").Value = 2888
$ws.Range("I40").Value = 2149.5
$ws.Range("J40").Value = 3183.4
$ws.Range("K40").Value = 2149.5
$ws.Range("L40").Value = 3183.4
$ws.Range("M40").Value = -2013.5
$ws.Range("N40").Value = -3455.4
$ws.Range("H61").Value = 1566.8334
$ws.Range("I61").Value = 1480.4
$ws.Range("K61").Value = 1480.4
$ws.Range("M61").Value = -1278.4
$ws.Range("H100").Value = 2213.25
$ws.Range("I100").Value = 1427.5
$ws.Range("K100").Value = 1427.5
$ws.Range("M100").Value = -886.5
$ws.Range("H113").Value = 1566.8334
$ws.Range("I113").Value = 1480.4
$ws.Range("K113").Value = 1480.4
$ws.Range("M113").Value = 689.5999999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3603
$ws.Range("I122").Value = 1762.1666
$ws.Range("J122").Value = 5812
$ws.Range("K122").Value = 5286.4998
$ws.Range("L122").Value = 17436
$ws.Range("M122").Value = -2836.4998
$ws.Range("N122").Value = -22336
$ws.Range("H126").Value = 1476
$ws.Range("I126").Value = 1202
$ws.Range("J126").Value = 1750
$ws.Range("K126").Value = 3606
$ws.Range("L126").Value = 5250
$ws.Range("M126").Value = -1136
$ws.Range("N126").Value = -10190

